# Refresh the crypto price/volume table (plus the swapped FTXToken/FraxShare
# and InjectiveProtocol/Cronos rows) per the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.778.33'
$ws.Range('E2').Value = '  +5.05%  '

$ws.Range('D3').Value = '2.279.47'
$ws.Range('E3').Value = '  +3.19%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '''234.25'
$ws.Range('E5').Value = '  +1.91%  '

$ws.Range('E6').Value = '  +3.55%  '

$ws.Range('D7').Value = '''65.50'
$ws.Range('E7').Value = '  +8.44%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').Value = '''0.430'
$ws.Range('E9').Value = '  +6.74%  '

$ws.Range('D10').Value = '''0.104'
$ws.Range('E10').Value = '  +17.16%  '

$ws.Range('D11').Value = '''57.68'
$ws.Range('E11').Value = '  +0.96%  '

$ws.Range('D12').Value = '''26.33'
$ws.Range('E12').Value = '  +18.60%  '

$ws.Range('D13').Value = '''0.103'
$ws.Range('E13').Value = '  +0.05%  '

$ws.Range('D14').Value = '2.617.80'
$ws.Range('E14').Value = '  +3.18%  '

$ws.Range('D15').Value = '''15.81'
$ws.Range('E15').Value = '  +2.54%  '

$ws.Range('D16').Value = '''5.98'
$ws.Range('E16').Value = '  +5.59%  '

$ws.Range('D17').Value = '''0.832'
$ws.Range('E17').Value = '  +4.72%  '

$ws.Range('D18').Value = '2.277.09'
$ws.Range('E18').Value = '  +3.26%  '

$ws.Range('D19').Value = '43.675.13'
$ws.Range('E19').Value = '  +4.74%  '

$ws.Range('D20').Value = '0.0₃0996'
$ws.Range('E20').Value = '  +10.58%  '

$ws.Range('D21').Value = '''74.43'
$ws.Range('E21').Value = '  +3.23%  '

$ws.Range('E22').Value = '  +1.72%  '

$ws.Range('D23').Value = '''262.72'
$ws.Range('E23').Value = '  +8.32%  '

$ws.Range('D25').Value = '''2.51'
$ws.Range('E25').Value = '  +6.73%  '

$ws.Range('D26').Value = '''2.28'
$ws.Range('E26').Value = '  -0.18%  '

$ws.Range('D27').Value = '''10.15'
$ws.Range('E27').Value = '  +5.00%  '

$ws.Range('D28').Value = '''172.66'
$ws.Range('E28').Value = '  +1.93%  '

$ws.Range('D29').Value = '''21.16'
$ws.Range('E29').Value = '  +7.00%  '

$ws.Range('E30').Value = '  -2.05%  '

$ws.Range('E31').Value = '  -0.45%  '

$ws.Range('E32').Value = '  +8.10%  '

$ws.Range('D33').Value = '''0.124'
$ws.Range('E33').Value = '  +2.87%  '

$ws.Range('D34').Value = '''0.0688'
$ws.Range('E34').Value = '  +5.97%  '

$ws.Range('D35').Value = '''5.10'
$ws.Range('E35').Value = '  +1.71%  '

$ws.Range('D36').Value = '''4.78'
$ws.Range('E36').Value = '  +3.18%  '

$ws.Range('D37').Value = '''6.81'
$ws.Range('E37').Value = '  +7.42%  '

$ws.Range('D38').Value = '''3.80'
$ws.Range('E38').Value = '  +7.40%  '

$ws.Range('E39').Value = '  -0.03%  '

$ws.Range('E40').Value = '  +4.40%  '

$ws.Range('E41').Value = '  +0.04%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''8.44'
$ws.Range('E42').Value = '  -1.37%  '

$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D43').Value = '''4.54'
$ws.Range('E43').Value = '  +3.19%  '

$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = '''0.0980'
$ws.Range('E44').Value = '  +2.56%  '

$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '''17.49'
$ws.Range('E45').Value = '  +6.62%  '

$ws.Range('D46').Value = '''10.44'
$ws.Range('E46').Value = '  +21.62%  '

$ws.Range('D47').Value = '''98.69'
$ws.Range('E47').Value = '  +1.45%  '

$ws.Range('E48').Value = '  +0.95%  '

$ws.Range('D49').Value = '1.479.90'
$ws.Range('E49').Value = '  +0.93%  '

$ws.Range('D50').Value = '''2.37'
$ws.Range('E50').Value = '  +7.37%  '

$ws.Range('D51').Value = '''0.000209'
$ws.Range('E51').Value = '  -12.75%  '
